# Update "想去人数" (F column) counts for a few events across sheets.
# Sheet 1: 展览 (Exhibitions)
# Sheet 2: 演出 (Performances)
# Sheet 3: 本地生活 (Local life)
# Sheet 4: 全部类型 (All types, aggregate of the above)

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F6").Value = 403   # 广州·炎焱动漫展 402 -> 403
$wsExhibit.Range("F11").Value = 6214 # 广州·AP动漫游戏嘉年华 6211 -> 6214
$wsExhibit.Range("F16").Value = 552  # 广州·运动番only 551 -> 552

# 演出 sheet updates
$wsShow.Range("F3").Value = 277      # 广州·KANAKO ITO&AYANE 2024 LIVE 276 -> 277

# 本地生活 sheet updates
$wsLocal.Range("F2").Value = 285     # 广州·NIJISANJI EN 官方授权主题店 284 -> 285

# 全部类型 sheet updates (aggregate duplicates of the above events)
$wsAll.Range("F2").Value = 285       # 广州·NIJISANJI EN 官方授权主题店 284 -> 285
$wsAll.Range("F8").Value = 403       # 广州·炎焱动漫展 402 -> 403
$wsAll.Range("F13").Value = 6214     # 广州·AP动漫游戏嘉年华 6211 -> 6214
$wsAll.Range("F16").Value = 277      # 广州·KANAKO ITO&AYANE 2024 LIVE 276 -> 277
$wsAll.Range("F19").Value = 552      # 广州·运动番only 551 -> 552
